# Update "想去人数" (interested-people count) figures on both the
# "展览" and "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 1619
    5  = 6145
    7  = 243
    12 = 258
    13 = 5614
    14 = 10320
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
